$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "set_impression_of_dialogue(&impression)" action-column cells are no
# longer needed now that the final utterance itself asks the LLM to generate the
# impression text, so remove those cells entirely.
$ws.Range("G20").Clear()
$ws.Range("G21").Clear()
$ws.Range("G22").Clear()

# Replace the old "{impression} . Thank you for your time." templated utterance
# with one that asks the model to generate the impression directly.
$ws.Range("C25").Value = '{$"Generate a short utterance to say the system''s impression."} Thank you for your time.'

# That row's text got longer, so Excel has given it an explicit row height.
$ws.Rows.Item(25).RowHeight = 45

# The active selection moved from F5 to C5.
$ws.Range("C5").Select()
